# data_required config: add columns to specify handling of NA values in calculations
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data_required")

# New header cells for the two added columns (AP, AQ) following the existing
# "A.*" / "B.*" naming convention used throughout row 1.
$ws.Range("AP1").Value = "A.value_na"
$ws.Range("AQ1").Value = "B.value_na"

# Body rows 2-56 and 59 default to the text "NA" (matches the rest of the
# table, which already uses "NA" as the default/placeholder string).
$ws.Range("AP2:AQ56").Value = "NA"
$ws.Range("AP59:AQ59").Value = "NA"

# Rows 57-58 (the "A+B" calculation rows) get a numeric 0 instead of "NA".
$ws.Range("AP57:AQ58").Value = 0

# Move the active selection onto the newly added column so the saved view
# reflects the edit (mirrors the author's workbookView/selection change).
$ws.Activate()
$ws.Range("AQ1").Select()
